$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$c = $ws.Range("A300:A305")
$c.Interior.Color = 16777215
